$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room above the Causeway row for two new Downtown showtimes (D02, D03)
# plus the relocated Downtown/D01 booking (row 7 slides down to row 10).
$ws.Rows("3:5").Insert()

# Pull the Downtown/D01 booking (now at row 10) up to row 5, preserving its
# original cell types (text "2", text date) via a straight range copy, then
# fix its time - it duplicated another Downtown/D01 showtime's slot.
$ws.Range("A10:E10").Copy($ws.Range("A5:E5"))
$ws.Range("E5").Value = 1500

# Remove the now-superseded original booking row.
$ws.Rows("10:10").Delete()

# The Causeway booking's movie number was stored as text - normalize it to a
# plain number like the rest of the "movie" column.
$ws.Range("A6").Value = 2

# New showtime: Downtown D02
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Downtown"
$ws.Range("C3").Value = "D02"
$ws.Range("D3").Value = 10122022
$ws.Range("E3").Value = 1500

# New showtime: Downtown D03
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Downtown"
$ws.Range("C4").Value = "D03"
$ws.Range("D4").Value = 10122022
$ws.Range("E4").Value = 1300

$ws.Range("C14").Select()
